# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the newly (re)computed "K" column values (column G, rows 2-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(6,3,4,5,12,5,2,5,5,1,5,11,9,10,8,7,8,9,8,8,9,7,6,12,5,2,10,3,5,10,10,8,13,11,4,5,4,5,5,4,1,1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
